$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Helper: replace the text of a whole paragraph range while avoiding the
# engine's "common prefix is kept as a separate run" behaviour. We first
# stomp the paragraph with a placeholder that shares no characters with
# either the old or the new text, then set the real text, so the result is
# always a single clean run.
# ---------------------------------------------------------------------------
function Set-ParagraphText($para, [string]$text) {
    $para.Text = "@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@"
    $para.Text = $text
}

# ===========================================================================
# Slide 2 ("Standards used in our project" -> "Assets used in the project")
# ===========================================================================
$s2 = $p.Slides.Item(2)

# Title
$title2 = $s2.Shapes.Item(1).TextFrame.TextRange
Set-ParagraphText $title2.Paragraphs(1,1) "Assets used in the project"

# Content placeholder
$body2 = $s2.Shapes.Item(2)

# Reset the "shrink text on overflow" cached scale (fontScale/lnSpcReduction)
# back to a plain <a:normAutofit/> now that the text is shorter.
$body2.TextFrame.AutoSize = 2

$tr2 = $body2.TextFrame.TextRange

# Paragraph 1: "Using GitHub to save our project so that the team can work together"
#           -> "Use of GitHub for team collaboration"
Set-ParagraphText $tr2.Paragraphs(1,1) "Use of GitHub for team collaboration"

# Paragraph 2: "Using GearHost to save our database on the cloud so that the
# team can <br/>have access to the same data"
#           -> "Database cloud-hosted on GearHost"
# Keep the existing "GearHost" run (with err="1") untouched; only touch the
# leading "Using " run and delete everything after "GearHost".
$para2 = $tr2.Paragraphs(2,1)
$len2 = $para2.Length
$trail2 = $para2.Characters(15, $len2 - 1 - 14)
$trail2.Text = ""
$para2 = $tr2.Paragraphs(2,1)
$lead2 = $para2.Characters(1, 6)
$lead2.Text = "Database cloud-hosted on "

# ===========================================================================
# Slide 3 ("Standards used in our project" -> "Assets used in the project (contd.)")
# ===========================================================================
$s3 = $p.Slides.Item(3)

# Title: also gains centre alignment
$title3Shape = $s3.Shapes.Item(1)
$title3 = $title3Shape.TextFrame.TextRange
Set-ParagraphText $title3.Paragraphs(1,1) "Assets used in the project (contd.)"
$title3.Paragraphs(1,1).ParagraphFormat.Alignment = 2

$tr3 = $s3.Shapes.Item(2).TextFrame.TextRange

# Paragraph 1: "Using JDBC API" -> "JDBC API for database connectivity and access"
Set-ParagraphText $tr3.Paragraphs(1,1) "JDBC API for database connectivity and access"

# Paragraph 3: "Using MySQL Database" -> "MySQL Database"
Set-ParagraphText $tr3.Paragraphs(3,1) "MySQL Database"

# ===========================================================================
# Slide 5 ("Lessons learned")
# ===========================================================================
$s5 = $p.Slides.Item(5)
$tr5 = $s5.Shapes.Item(2).TextFrame.TextRange

Set-ParagraphText $tr5.Paragraphs(1,1) "Implementing MVC architecture"
Set-ParagraphText $tr5.Paragraphs(2,1) "Login/Logout and user role functionality"
Set-ParagraphText $tr5.Paragraphs(3,1) "User sessions"
Set-ParagraphText $tr5.Paragraphs(4,1) "Implementing a search function"

# Paragraph 5: "How " + "to use " + "include for consistent layout"
#           -> "Using includes to keep " + "layout consistent"
$para5 = $tr5.Paragraphs(5,1)
$len5 = $para5.Length
$lead5 = $para5.Characters(1, 4)
$lead5.Text = "Using includes to keep "
$para5 = $tr5.Paragraphs(5,1)
$len5 = $para5.Length
$tail5 = $para5.Characters(25, $len5 - 1 - 24)
$tail5.Text = "layout consistent"
